$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.448.25'
$ws.Range("E2").Value = '  -0.96%  '

$ws.Range("D3").Value = '3.534.75'
$ws.Range("E3").Value = '  -2.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '196.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '582.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.614'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.20%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.204'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.627'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '51.91'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000287'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.94%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.84%  '

$ws.Range("B14").Value = 'BitcoinCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '682.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +15.25%  '

$ws.Range("D15").Value = '4.103.02'
$ws.Range("E15").Value = '  -2.33%  '

$ws.Range("D16").Value = '69.508.77'
$ws.Range("E16").Value = '  -1.19%  '

$ws.Range("D17").Value = '3.560.75'
$ws.Range("E17").Value = '  -1.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.46'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.55'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.121'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.968'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '18.20'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '107.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.64%  '

$ws.Range("E24").Value = '  +1.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.42%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.72%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.12%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.38'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.81%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.111'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.17%  '

$ws.Range("D35").Value = '3.794.62'
$ws.Range("E35").Value = '  -3.66%  '

$ws.Range("B36").Value = 'Stacks'
$ws.Range("C36").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.09%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0813'
$ws.Range("E37").Value = '  -9.51%  '

$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.95'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '498.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.373'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.135'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '34.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.01%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0462'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.14%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.38'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.137'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.81%  '

$ws.Range("E48").Value = '  -0.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.83%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +19.96%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +60.05%  '
